$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.006.30'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '1.704.20'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4002'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.32%  '
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.474'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08823'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.10'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.491'
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.978'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.59%  '
$ws.Range("D17").Value = '1.675.22'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.323'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").Value = '25.009.94'
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.402'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.941'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.60'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.059'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '152.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.420'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.693'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +21.62%  '
$ws.Range("D33").Value = '1.886.57'
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08632'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.69%  '
$ws.Range("E35").Value = '  +3.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.051'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.200'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2921'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.91%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09693'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8251'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.05'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.484'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.694'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7379'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.09192'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.252'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.406'
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.99%  '
